$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Updated TPM-derived values for rows 2-10 (columns E through T)
# Row 2
$ws.Range("E2").Value = 1
$ws.Range("F2").Value = 0.3333333333333333
$ws.Range("G2").Value = 0.02185066666666667
$ws.Range("H2").Value = 0.065552
$ws.Range("I2").Value = 0.02597345993572409
$ws.Range("J2").Value = 0.02597345993572409
$ws.Range("K2").Value = 2
$ws.Range("L2").Value = 0.6666666666666666
$ws.Range("M2").Value = 0.20543
$ws.Range("N2").Value = 0.61629
$ws.Range("O2").Value = 0.0348838848157659
$ws.Range("P2").Value = 0.0348838848157659
$ws.Range("Q2").Value = 0.004488782453333334
$ws.Range("R2").Value = 0.04039904208
$ws.Range("S2").Value = 0.0009060551846647094
$ws.Range("T2").Value = 0.0009060551846647094

# Row 3
$ws.Range("E3").Value = 1
$ws.Range("F3").Value = 0.3333333333333333
$ws.Range("G3").Value = 0.02185066666666667
$ws.Range("H3").Value = 0.065552
$ws.Range("I3").Value = 0.02597345993572409
$ws.Range("J3").Value = 0.02597345993572409
$ws.Range("O3").Value = 0.01910092077856117
$ws.Range("P3").Value = 0.01910092077856117
$ws.Range("Q3").Value = 0.002457864956444445
$ws.Range("R3").Value = 0.022120784608
$ws.Range("S3").Value = 0.0004961170005773984
$ws.Range("T3").Value = 0.0004961170005773984

# Row 4
$ws.Range("E4").Value = 1
$ws.Range("F4").Value = 0.3333333333333333
$ws.Range("G4").Value = 0.02185066666666667
$ws.Range("H4").Value = 0.065552
$ws.Range("I4").Value = 0.02597345993572409
$ws.Range("J4").Value = 0.02597345993572409
$ws.Range("O4").Value = 0.9460151944056729
$ws.Range("P4").Value = 0.9460151944056729
$ws.Range("Q4").Value = 0.121731178384
$ws.Range("R4").Value = 1.095580605456
$ws.Range("S4").Value = 0.02457128775048198
$ws.Range("T4").Value = 0.02457128775048198

# Row 5
$ws.Range("I5").Value = 0.6906391812052189
$ws.Range("J5").Value = 0.6906391812052189
$ws.Range("K5").Value = 2
$ws.Range("L5").Value = 0.6666666666666666
$ws.Range("M5").Value = 0.20543
$ws.Range("N5").Value = 0.61629
$ws.Range("O5").Value = 0.0348838848157659
$ws.Range("P5").Value = 0.0348838848157659
$ws.Range("Q5").Value = 0.1193575690666667
$ws.Range("R5").Value = 1.0742181216
$ws.Range("S5").Value = 0.02409217764641773
$ws.Range("T5").Value = 0.02409217764641773

# Row 6
$ws.Range("I6").Value = 0.6906391812052189
$ws.Range("J6").Value = 0.6906391812052189
$ws.Range("O6").Value = 0.01910092077856117
$ws.Range("P6").Value = 0.01910092077856117
$ws.Range("S6").Value = 0.01319184428677124
$ws.Range("T6").Value = 0.01319184428677124

# Row 7
$ws.Range("I7").Value = 0.6906391812052189
$ws.Range("J7").Value = 0.6906391812052189
$ws.Range("O7").Value = 0.9460151944056729
$ws.Range("P7").Value = 0.9460151944056729
$ws.Range("S7").Value = 0.65335515927203
$ws.Range("T7").Value = 0.65335515927203

# Row 8
$ws.Range("I8").Value = 0.283387358859057
$ws.Range("J8").Value = 0.283387358859057
$ws.Range("K8").Value = 2
$ws.Range("L8").Value = 0.6666666666666666
$ws.Range("M8").Value = 0.20543
$ws.Range("N8").Value = 0.61629
$ws.Range("O8").Value = 0.0348838848157659
$ws.Range("P8").Value = 0.0348838848157659
$ws.Range("Q8").Value = 0.04897553915
$ws.Range("R8").Value = 0.4407798523499999
$ws.Range("S8").Value = 0.009885651984683459
$ws.Range("T8").Value = 0.009885651984683459

# Row 9
$ws.Range("I9").Value = 0.283387358859057
$ws.Range("J9").Value = 0.283387358859057
$ws.Range("O9").Value = 0.01910092077856117
$ws.Range("P9").Value = 0.01910092077856117
$ws.Range("S9").Value = 0.005412959491212533
$ws.Range("T9").Value = 0.005412959491212533

# Row 10
$ws.Range("I10").Value = 0.283387358859057
$ws.Range("J10").Value = 0.283387358859057
$ws.Range("O10").Value = 0.9460151944056729
$ws.Range("P10").Value = 0.9460151944056729
$ws.Range("S10").Value = 0.268088747383161
$ws.Range("T10").Value = 0.268088747383161
